$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update site_id for row 2 from "M1" to "M0"
$ws.Range("A2").Value = "M0"

# Update the selection to match the saved view state
$ws.Range("B10").Select()
